$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 900
$ws.Range("J45").Value = 900
$ws.Range("L45").Value = 2700
$ws.Range("N45").Value = -3084
$ws.Range("H129").Value = 186228.3
$ws.Range("I129").Value = 248.5
$ws.Range("J129").Value = 193381.36
$ws.Range("K129").Value = 745.5
$ws.Range("L129").Value = 580144.08
$ws.Range("M129").Value = 4254.5
$ws.Range("N129").Value = -590144.08
$ws.Range("H132").Value = 2262.617
$ws.Range("I132").Value = 2513.7693
$ws.Range("J132").Value = 1038.25
$ws.Range("K132").Value = 7541.3079
$ws.Range("L132").Value = 3114.75
$ws.Range("M132").Value = -5011.3079
$ws.Range("N132").Value = -8174.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2393
$ws.Range("I45").Value = 1784.4615
$ws.Range("K45").Value = 1784.4615
$ws.Range("M45").Value = -1407.4615
$ws.Range("H102").Value = 1137
$ws.Range("I102").Value = 1079.125
$ws.Range("J102").Value = 1600
$ws.Range("K102").Value = 1079.125
$ws.Range("L102").Value = 1600
$ws.Range("M102").Value = 542.875
$ws.Range("N102").Value = -4844
$ws.Range("H110").Value = 619.6667
$ws.Range("I110").Value = 554.75
$ws.Range("J110").Value = 749.5
$ws.Range("K110").Value = 554.75
$ws.Range("L110").Value = 749.5
$ws.Range("M110").Value = 1490.25
$ws.Range("N110").Value = -4839.5
$ws.Range("H122").Value = 1397.326
$ws.Range("I122").Value = 1242.8422
$ws.Range("K122").Value = 3728.5266
$ws.Range("M122").Value = -1278.5266
$ws.Range("H132").Value = 16618.117
$ws.Range("I132").Value = 1582.2858
$ws.Range("J132").Value = 86785.336
$ws.Range("K132").Value = 4746.857400000001
$ws.Range("L132").Value = 260356.008
$ws.Range("M132").Value = -2216.857400000001
$ws.Range("N132").Value = -265416.008

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1013.625
$ws.Range("I94").Value = 851.5
$ws.Range("K94").Value = 851.5
$ws.Range("M94").Value = -400.5
$ws.Range("H134").Value = 2705.868
$ws.Range("I134").Value = 2621.0833
$ws.Range("J134").Value = 3519.8
$ws.Range("K134").Value = 7863.249899999999
$ws.Range("L134").Value = 10559.4
$ws.Range("M134").Value = -5328.249899999999
$ws.Range("N134").Value = -15629.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 551
$ws.Range("I10").Value = 551
$ws.Range("K10").Value = 551
$ws.Range("M10").Value = -412
$ws.Range("H31").Value = 3438.8096
$ws.Range("I31").Value = 3693.5
$ws.Range("J31").Value = 3359.2188
$ws.Range("K31").Value = 3693.5
$ws.Range("L31").Value = 3359.2188
$ws.Range("M31").Value = -3398.5
$ws.Range("N31").Value = -3949.2188
$ws.Range("H34").Value = 3438.8096
$ws.Range("I34").Value = 3693.5
$ws.Range("J34").Value = 3359.2188
$ws.Range("K34").Value = 3693.5
$ws.Range("L34").Value = 3359.2188
$ws.Range("M34").Value = -3491.5
$ws.Range("N34").Value = -3763.2188
$ws.Range("H58").Value = 18883
$ws.Range("J58").Value = 34994.2
$ws.Range("L58").Value = 34994.2
$ws.Range("N58").Value = -35400.2
$ws.Range("H136").Value = 18883
$ws.Range("J136").Value = 34994.2
$ws.Range("L136").Value = 104982.6
$ws.Range("N136").Value = -110082.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1283.7778
$ws.Range("I5").Value = 828.2222
$ws.Range("J5").Value = 1739.3334
$ws.Range("K5").Value = 2484.6666
$ws.Range("L5").Value = 5218.0002
$ws.Range("M5").Value = -2372.6666
$ws.Range("N5").Value = -5442.0002
$ws.Range("H75").Value = 1957.5
$ws.Range("J75").Value = 1915
$ws.Range("L75").Value = 5745
$ws.Range("N75").Value = -7741
$ws.Range("H78").Value = 1957.5
$ws.Range("J78").Value = 1915
$ws.Range("L78").Value = 17235
$ws.Range("N78").Value = -27219
$ws.Range("H117").Value = 1198.4445
$ws.Range("J117").Value = 1151.4
$ws.Range("L117").Value = 3454.2
$ws.Range("N117").Value = -10338.2
$ws.Range("H131").Value = 754.24
$ws.Range("I131").Value = 650
$ws.Range("J131").Value = 755.2929
$ws.Range("K131").Value = 1950
$ws.Range("L131").Value = 2265.8787
$ws.Range("M131").Value = 3090
$ws.Range("N131").Value = -12345.8787
$ws.Range("H135").Value = 1283.7778
$ws.Range("I135").Value = 828.2222
$ws.Range("J135").Value = 1739.3334
$ws.Range("K135").Value = 7453.999800000001
$ws.Range("L135").Value = 15654.0006
$ws.Range("M135").Value = -4918.999800000001
$ws.Range("N135").Value = -20724.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3393.4167
$ws.Range("I80").Value = 2699.25
$ws.Range("J80").Value = 3740.5
$ws.Range("K80").Value = 2699.25
$ws.Range("L80").Value = 3740.5
$ws.Range("M80").Value = -1701.25
$ws.Range("N80").Value = -5736.5
$ws.Range("H83").Value = 3393.4167
$ws.Range("I83").Value = 2699.25
$ws.Range("J83").Value = 3740.5
$ws.Range("K83").Value = 13496.25
$ws.Range("L83").Value = 18702.5
$ws.Range("M83").Value = -8504.25
$ws.Range("N83").Value = -28686.5
$ws.Range("H93").Value = 13121.429
$ws.Range("J93").Value = 13121.429
$ws.Range("L93").Value = 13121.429
$ws.Range("N93").Value = -16865.429
$ws.Range("H102").Value = 14287737
$ws.Range("I102").Value = 16668517
$ws.Range("K102").Value = 16668517
$ws.Range("M102").Value = -16666895
$ws.Range("H113").Value = 12803
$ws.Range("I113").Value = 21777.75
$ws.Range("J113").Value = 3828.25
$ws.Range("K113").Value = 21777.75
$ws.Range("L113").Value = 3828.25
$ws.Range("M113").Value = -19607.75
$ws.Range("N113").Value = -8168.25
$ws.Range("H122").Value = 45978036
$ws.Range("I122").Value = 15873990
$ws.Range("J122").Value = 125001160
$ws.Range("K122").Value = 47621970
$ws.Range("L122").Value = 375003480
$ws.Range("M122").Value = -47619520
$ws.Range("N122").Value = -375008380
$ws.Range("H123").Value = 5126.6
$ws.Range("J123").Value = 20326
$ws.Range("L123").Value = 20326
$ws.Range("N123").Value = -25226
$ws.Range("H126").Value = 5649.04
$ws.Range("I126").Value = 4491.1763
$ws.Range("J126").Value = 8109.5
$ws.Range("K126").Value = 13473.5289
$ws.Range("L126").Value = 24328.5
$ws.Range("M126").Value = -11003.5289
$ws.Range("N126").Value = -29268.5
$ws.Range("H132").Value = 22989.166
$ws.Range("I132").Value = 1952.05
$ws.Range("J132").Value = 128174.75
$ws.Range("K132").Value = 5856.15
$ws.Range("L132").Value = 384524.25
$ws.Range("M132").Value = -3326.15
$ws.Range("N132").Value = -389584.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4181.727
$ws.Range("I7").Value = 3700
$ws.Range("K7").Value = 3700
$ws.Range("M7").Value = -3588
$ws.Range("H40").Value = 2763.1143
$ws.Range("I40").Value = 2269.8076
$ws.Range("J40").Value = 4188.222
$ws.Range("K40").Value = 2269.8076
$ws.Range("L40").Value = 4188.222
$ws.Range("M40").Value = -2133.8076
$ws.Range("N40").Value = -4460.222
$ws.Range("H68").Value = 2328.7
$ws.Range("I68").Value = 2032.6666
$ws.Range("J68").Value = 2455.5715
$ws.Range("K68").Value = 2032.6666
$ws.Range("L68").Value = 2455.5715
$ws.Range("M68").Value = -1283.6666
$ws.Range("N68").Value = -3953.5715
$ws.Range("H71").Value = 2328.7
$ws.Range("I71").Value = 2032.6666
$ws.Range("J71").Value = 2455.5715
$ws.Range("K71").Value = 10163.333
$ws.Range("L71").Value = 12277.8575
$ws.Range("M71").Value = -6419.333000000001
$ws.Range("N71").Value = -19765.8575
$ws.Range("H126").Value = 4181.727
$ws.Range("I126").Value = 3700
$ws.Range("K126").Value = 11100
$ws.Range("M126").Value = -8630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3599.8
$ws.Range("I62").Value = 3249.5
$ws.Range("J62").Value = 3833.3333
$ws.Range("K62").Value = 3249.5
$ws.Range("L62").Value = 3833.3333
$ws.Range("M62").Value = -2625.5
$ws.Range("N62").Value = -5081.3333
$ws.Range("H65").Value = 3599.8
$ws.Range("I65").Value = 3249.5
$ws.Range("J65").Value = 3833.3333
$ws.Range("K65").Value = 16247.5
$ws.Range("L65").Value = 19166.6665
$ws.Range("M65").Value = -13127.5
$ws.Range("N65").Value = -25406.6665
$ws.Range("H132").Value = 1362.0588
$ws.Range("I132").Value = 846.3333
$ws.Range("K132").Value = 2538.9999
$ws.Range("M132").Value = -8.999899999999798
